$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 4.5
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62
$ws.Range("AR3").Value = 51

# Row 4
$ws.Range("A4").Value = 'G6BcHP1l'
$ws.Range("C4").Value = '22:30'
$ws.Range("D4").Value = 'COLOMBIA - PRIMERA A'
$ws.Range("E4").Value = 'Millonarios'
$ws.Range("F4").Value = 'Pereira'
$ws.Range("G4").Value = 1.57
$ws.Range("H4").Value = 3.5
$ws.Range("I4").Value = 6.25
$ws.Range("J4").Value = 2.25
$ws.Range("L4").Value = 6.5
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 7.5
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 2.75
$ws.Range("Q4").Value = 2.25
$ws.Range("R4").Value = 1.62
$ws.Range("S4").Value = 1.5
$ws.Range("T4").Value = 2.5
$ws.Range("U4").Value = 2.25
$ws.Range("V4").Value = 1.57
$ws.Range("W4").Value = 5.5
$ws.Range("X4").Value = 6.5
$ws.Range("Y4").Value = 9
$ws.Range("Z4").Value = 11
$ws.Range("AA4").Value = 15
$ws.Range("AB4").Value = 34
$ws.Range("AC4").Value = 7.5
$ws.Range("AD4").Value = 7
$ws.Range("AE4").Value = 21
$ws.Range("AF4").Value = 81
$ws.Range("AH4").Value = 13
$ws.Range("AI4").Value = 29
$ws.Range("AJ4").Value = 21
$ws.Range("AK4").Value = 67
$ws.Range("AL4").Value = 51
$ws.Range("AM4").Value = 51
$ws.Range("AN4").Value = 3.4
$ws.Range("AO4").Value = 8.5
$ws.Range("AP4").Value = 23
$ws.Range("AQ4").Value = 29
$ws.Range("AS4").Value = 201
$ws.Range("AT4").Value = 2.5
$ws.Range("AU4").Value = 9.5
$ws.Range("AV4").Value = 81
$ws.Range("AW4").Value = 7
$ws.Range("AX4").Value = 34
$ws.Range("AY4").Value = 41
$ws.Range("AZ4").Value = 151
$ws.Range("BA4").Value = 201
$ws.Range("BB4").Value = 501
$ws.Range("BC4").Value = 126

# Row 6
$ws.Range("G6").Value = 2.2
$ws.Range("I6").Value = 3.1
$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 2.75
$ws.Range("Q6").Value = 2.25
$ws.Range("R6").Value = 1.62
$ws.Range("W6").Value = 7
$ws.Range("X6").Value = 10
$ws.Range("AA6").Value = 21
$ws.Range("AN6").Value = 4.33

# Row 7
$ws.Range("H7").Value = 3.65
$ws.Range("J7").Value = 2.15
$ws.Range("L7").Value = 5.4
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 7.1
$ws.Range("O7").Value = 1.32
$ws.Range("P7").Value = 3.1
$ws.Range("Q7").Value = 1.95
$ws.Range("U7").Value = 1.95
$ws.Range("W7").Value = 6.2
$ws.Range("Z7").Value = 11.5
$ws.Range("AC7").Value = 7.1
$ws.Range("AD7").Value = 7.2
$ws.Range("AJ7").Value = 17.5
$ws.Range("AL7").Value = 60
$ws.Range("AM7").Value = 65
$ws.Range("AN7").Value = 3.35
$ws.Range("AO7").Value = 7.7
$ws.Range("AP7").Value = 18
$ws.Range("AQ7").Value = 25
$ws.Range("AU7").Value = 8
$ws.Range("AX7").Value = 32
$ws.Range("AY7").Value = 37
